# Add 2022-Q4 data.
#
# before.xlsx has two sheets: "总计" (totals) and "2022-Q3" (the only
# quarterly detail sheet so far). This change:
#   1. Turns the existing "2022-Q3" detail sheet into "2022-Q4"
#      (new holdings data), re-using its sheetId/rId (=> sheetId 2).
#   2. Adds a brand new sheet named "2022-Q3" right after it, populated
#      with the data the "2022-Q3" sheet used to hold (=> sheetId 3,
#      appended last so sheetId numbering matches).
#   3. Appends a matching "2022-Q3" row to the "总计" sheet (which
#      already got its row 2 relabeled "2022-Q4").

$wb = $excel.ActiveWorkbook
$sheetTotal = $wb.Worksheets.Item(1)   # "总计"
$sheetQ = $wb.Worksheets.Item(2)       # currently "2022-Q3"

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet (after the current Q3/soon-Q4 sheet)
#    and clone the OLD Q3 data + formatting into it before we overwrite
#    the source sheet's content.
# ---------------------------------------------------------------------
$newQ3 = $wb.Worksheets.Add($null, $sheetQ)

$sheetQ.Range("B1:H2").Copy()
$newQ3.Range("B1:H2").PasteSpecial(-4122)
$newQ3.Range("B1:H2").PasteSpecial(-4163)

$sheetQ.Range("A2").Copy()
$newQ3.Range("A2").PasteSpecial(-4122)
$newQ3.Range("A2").PasteSpecial(-4163)

# ---------------------------------------------------------------------
# 2) Rename the ORIGINAL sheet to "2022-Q4" first (so the rename below
#    doesn't collide with the still-existing "2022-Q3" name), then
#    rename the clone to "2022-Q3".
# ---------------------------------------------------------------------
$sheetQ.Name = "2022-Q4"
$newQ3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 3) Re-style the header row + A2 of the (now) "2022-Q4" sheet to match
#    the "总计" sheet's style (matches the authored workbook), then
#    write the new Q4 values.
# ---------------------------------------------------------------------
$sheetTotal.Range("B1:D1").Copy()
$sheetQ.Range("B1:H1").PasteSpecial(-4122)
$sheetTotal.Range("A2").Copy()
$sheetQ.Range("A2").PasteSpecial(-4122)

$sheetQ.Range("B1").Value = "基金代码"
$sheetQ.Range("C1").Value = "基金名称"
$sheetQ.Range("D1").Value = "基金规模"
$sheetQ.Range("E1").Value = "股票总仓位"
$sheetQ.Range("F1").Value = "仓位占比"
$sheetQ.Range("G1").Value = "持有市值(亿元)"
$sheetQ.Range("H1").Value = "仓位排名"

$sheetQ.Range("A2").Value = 0
# Keep these as literal text (not auto-converted numbers) by assigning
# a string-literal formula, then freezing it to a plain value.
$sheetQ.Range("B2").Formula = "=""002567"""
$sheetQ.Range("C2").Value = "大成国家安全主题灵活配置混合"
$sheetQ.Range("D2").Formula = "=""0.41"""
$sheetQ.Range("E2").Formula = "=""54.37"""
$sheetQ.Range("F2").Formula = "=""3.56"""
$sheetQ.Range("G2").Formula = "=""0.0146"""
$sheetQ.Range("B2:G2").Copy()
$sheetQ.Range("B2:G2").PasteSpecial(-4163)
$sheetQ.Range("H2").Value = 8

# ---------------------------------------------------------------------
# 4) Update "总计": row 2 becomes "2022-Q4" (data unchanged), and a new
#    row 3 repeats the same counts labelled "2022-Q3".
# ---------------------------------------------------------------------
$sheetTotal.Range("B2").Value = "2022-Q4"

$sheetTotal.Range("A2").Copy()
$sheetTotal.Range("A3").PasteSpecial(-4122)
$sheetTotal.Range("A3").Value = 1
$sheetTotal.Range("B3").Value = "2022-Q3"
$sheetTotal.Range("C3").Value = 1
$sheetTotal.Range("D3").Value = 0.01

Write-Host "done"
